# Automatic update of files.
# - Column C ("Förändrad") bumped from 46070 -> 46072 for every data row (2..25).
# - Rows 11..19 got re-sorted: their Beteckning (A), Datum (B) and Area (G)
#   values now come from a different source row (see mapping below).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Bump "Förändrad" (column C) for all data rows from 46070 to 46072.
for ($row = 2; $row -le 25; $row++) {
    $ws.Cells.Item($row, 3).Value2 = 46072
}

# 2) Rows 11..19 are re-ordered. Capture the current (pre-sort) A/B/G values
#    for these rows first, then write them back out in their new positions.
$colA = @{}
$colB = @{}
$colG = @{}
for ($row = 11; $row -le 19; $row++) {
    $colA[$row] = $ws.Cells.Item($row, 1).Value2
    $colB[$row] = $ws.Cells.Item($row, 2).Value2
    $colG[$row] = $ws.Cells.Item($row, 7).Value2
}

# new row -> old row the data is taken from
$rowMap = @{
    11 = 16
    12 = 15
    13 = 17
    14 = 11
    15 = 12
    16 = 13
    17 = 18
    18 = 19
    19 = 14
}

foreach ($newRow in 11..19) {
    $srcRow = $rowMap[$newRow]
    $ws.Cells.Item($newRow, 1).Value = $colA[$srcRow]
    $ws.Cells.Item($newRow, 2).Value2 = $colB[$srcRow]
    $ws.Cells.Item($newRow, 7).Value2 = $colG[$srcRow]
}
